$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed/Updated) date column C for rows 2-101
# from 2023-09-01 (45170) to 2023-09-05 (45174), preserving the existing
# date formatting already applied to those cells.
$ws.Range("C2:C101").Value = 45174
